$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Layer0")
$ws1.Range("B2").Value = -0.5578981002433783
$ws1.Range("C2").Value = -0.7723424526762497
$ws1.Range("B3").Value = 0.3554350327846399
$ws1.Range("C3").Value = 0.2212628635966589
$ws1.Range("B4").Value = 0.8683777379430004
$ws1.Range("C4").Value = 0.8184912724751707

$ws2 = $wb.Worksheets.Item("Layer1")
$ws2.Range("B2").Value = -1.377181892596393
$ws2.Range("C2").Value = -0.3476465013015807
$ws2.Range("B3").Value = 0.5458976361397698
$ws2.Range("C3").Value = 0.9755672625228474
$ws2.Range("B4").Value = 0.3839077335506901
$ws2.Range("C4").Value = -0.3805352388642115
